# Insert a new weekly record at the top of the Coliflor price table (row 113),
# pushing the existing records (old rows 113-170) down to rows 114-171.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 113:170 down to 114:171, duplicating formatting
# (date style, etc.) from the row above - matches dimension A1:R170 -> A1:R171.
$ws.Rows(113).Insert()

# Populate the newly inserted row 113 with the new observation.
$ws.Cells.Item(113, 1).Value = 5
$ws.Cells.Item(113, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(113, 3).Value = "Maule"
$ws.Cells.Item(113, 4).Value = 44518
$ws.Cells.Item(113, 5).Value = 7
$ws.Cells.Item(113, 6).Value = 100112008
$ws.Cells.Item(113, 7).Value = "Coliflor"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 4000
$ws.Cells.Item(113, 11).Value = 500
$ws.Cells.Item(113, 12).Value = 500
$ws.Cells.Item(113, 13).Value = 500
$ws.Cells.Item(113, 14).Value = "$/unidad"
$ws.Cells.Item(113, 15).Value = "Región del Maule"
$ws.Cells.Item(113, 16).Value = 500
$ws.Cells.Item(113, 17).Value = 1
$ws.Cells.Item(113, 18).Value = "Hortaliza"
